# Add five additional "FTNC_DemandN" result sheets (debug sheets for
# checking database / df output), each cloning the layout/formatting of
# the original FTNC_Demand5 sheet with its own set of summary numbers.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$newSheetData = @(
    @{ Name = "FTNC_Demand51"; B2 = 12.14893524930748;  C2 = 184.5303067413828;  D2 = 0; E2 = 0;                 F2 = 196.6792419906901 },
    @{ Name = "FTNC_Demand52"; B2 = 14.7297870498615;   C2 = 189.4586169906902;  D2 = 0; E2 = 75.18848069396029;  F2 = 279.3768847345121 },
    @{ Name = "FTNC_Demand53"; B2 = 11.96169494459834;  C2 = 184.5360338881972;  D2 = 0; E2 = 0;                 F2 = 196.4977288327954 },
    @{ Name = "FTNC_Demand54"; B2 = 10.71282029085873;  C2 = 182.7747665751778;  D2 = 0; E2 = 0;                 F2 = 193.4875868660364 },
    @{ Name = "FTNC_Demand55"; B2 = 14.72926765927978;  C2 = 187.0297942760089;  D2 = 0; E2 = 8.28313049376839;   F2 = 210.0421924290569 }
)

foreach ($sheetInfo in $newSheetData) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $sheetInfo.Name

    # Clone the header row + row-label cell formatting from the original sheet.
    $ws1.Range("B1:F1").Copy()
    $ws.Range("B1").PasteSpecial(-4122)

    $ws1.Range("A2").Copy()
    $ws.Range("A2").PasteSpecial(-4122)

    $ws.Range("B1").Value = "In-vehicle"
    $ws.Range("C1").Value = "At-stop"
    $ws.Range("D1").Value = "Extra"
    $ws.Range("E1").Value = "Tardiness"
    $ws.Range("F1").Value = "Total"

    $ws.Range("A2").Value = "FTNC"
    $ws.Range("B2").Value = $sheetInfo.B2
    $ws.Range("C2").Value = $sheetInfo.C2
    $ws.Range("D2").Value = $sheetInfo.D2
    $ws.Range("E2").Value = $sheetInfo.E2
    $ws.Range("F2").Value = $sheetInfo.F2
}
